# Deploy the implementation guide:
#  - bump the "Date" metadata value on the Metadata sheet
#  - append three new concept rows to the "Concepts" sheet, matching the
#    existing look & feel (copy format + value from the last data row,
#    then overwrite Code/Display)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2023-03-21T11:43:18+00:00"

$concepts = $wb.Worksheets.Item("Concepts")

$newRows = @(
    @("TRATU", "Transcriptome"),
    @("TUHEM", "Hematological malignancies predisposition"),
    @("TUPED", "Pediatric cancer predisposition")
)

$lastRow = 11
$row = $lastRow + 1
foreach ($entry in $newRows) {
    $src = $concepts.Range("A" + $lastRow + ":D" + $lastRow)
    $dst = $concepts.Range("A" + $row + ":D" + $row)

    # Copy formatting (borders/fill/alignment) from the template row.
    $src.Copy()
    $dst.PasteSpecial(-4122)

    # Copy values too (keeps column A's "1" a text value, like the rest
    # of the table) then overwrite Code / Display with the new data.
    $src.Copy()
    $dst.PasteSpecial(-4163)

    $concepts.Cells.Item($row, 2).Value = $entry[0]
    $concepts.Cells.Item($row, 3).Value = $entry[1]

    $row = $row + 1
}
